$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item(1)

# New row 59 mirrors row 58 (same model run values), but represents a new
# baseline run "Baseline 2010-18 C174" (added after "Baseline 2010-18 C173").
$row = 59

$ws.Cells.Item($row, 1).Value() = "CW3M"
$ws.Cells.Item($row, 2).Value() = "Baseline 2010-18 C174"
$ws.Cells.Item($row, 3).Value() = "2010-18"

$ws.Cells.Item($row, 4).Value() = 1080.801350777778
$ws.Cells.Item($row, 5).Value() = 1901.5157334444443
$ws.Cells.Item($row, 6).Value() = 0.97970299999999988
$ws.Cells.Item($row, 7).Value() = 280.33542888888883
$ws.Cells.Item($row, 8).Value() = 9.775355222222224
$ws.Cells.Item($row, 9).Value() = 5.7424886666666666
$ws.Cells.Item($row, 10).Value() = 8.145128999999999
$ws.Cells.Item($row, 11).Value() = 645.86557344444441
$ws.Cells.Item($row, 12).Value() = 83.47062044444445
$ws.Cells.Item($row, 13).Value() = 1430.3511555555554
$ws.Cells.Item($row, 14).Value() = 1110.4303452222223
$ws.Cells.Item($row, 15).Value() = 4638.5936415555561
$ws.Cells.Item($row, 16).Value() = 27227.338324888889
$ws.Cells.Item($row, 17).Value() = -0.8872363333333334
$ws.Cells.Item($row, 18).Value() = -0.00025788888888888888
$ws.Cells.Item($row, 19).Value() = "2010-18"

# Match number formats from the row above (D:N "0.00", O:P "0", Q "0.00",
# R "0.000000"); columns A-C and S stay plain text/general. Note: unlike
# M58 (highlighted yellow), M59 is left with the plain "0.00" format/no fill.
$ws.Range("D" + $row + ":N" + $row).NumberFormat = "0.00"
$ws.Range("O" + $row + ":P" + $row).NumberFormat = "0"
$ws.Range("Q" + $row).NumberFormat = "0.00"
$ws.Range("R" + $row).NumberFormat = "0.000000"

$ws.Application.Goto($ws.Range("S" + $row))
